$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1576.1538
$ws.Range("I19").Value = 1618.8182
$ws.Range("K19").Value = 1618.8182
$ws.Range("M19").Value = -1443.8182
$ws.Range("H103").Value = 2213.7144
$ws.Range("I103").Value = 1832.6666
$ws.Range("J103").Value = 2499.5
$ws.Range("K103").Value = 5497.9998
$ws.Range("L103").Value = 7498.5
$ws.Range("M103").Value = -4911.9998
$ws.Range("N103").Value = -8670.5
$ws.Range("H104").Value = 1298.25
$ws.Range("J104").Value = 1945
$ws.Range("L104").Value = 5835
$ws.Range("N104").Value = -9329
$ws.Range("H112").Value = 1854701.8
$ws.Range("I112").Value = 1766.6666
$ws.Range("J112").Value = 2225288.8
$ws.Range("K112").Value = 5299.9998
$ws.Range("L112").Value = 6675866.399999999
$ws.Range("M112").Value = -4191.9998
$ws.Range("N112").Value = -6678082.399999999
$ws.Range("H113").Value = 6231.7827
$ws.Range("I113").Value = 8528.866
$ws.Range("J113").Value = 1924.75
$ws.Range("K113").Value = 8528.866
$ws.Range("L113").Value = 1924.75
$ws.Range("M113").Value = -5274.866
$ws.Range("N113").Value = -8432.75
$ws.Range("H131").Value = 4682.5713
$ws.Range("I131").Value = 490
$ws.Range("J131").Value = 6359.6
$ws.Range("K131").Value = 1470
$ws.Range("L131").Value = 19078.8
$ws.Range("M131").Value = 3570
$ws.Range("N131").Value = -29158.8
$ws.Range("H132").Value = 1659.7037
$ws.Range("I132").Value = 1710.3556
$ws.Range("K132").Value = 5131.066800000001
$ws.Range("M132").Value = -2601.066800000001
$ws.Range("H135").Value = 2730.7222
$ws.Range("I135").Value = 1596.5834
$ws.Range("J135").Value = 4999
$ws.Range("K135").Value = 14369.2506
$ws.Range("L135").Value = 44991
$ws.Range("M135").Value = -11834.2506
$ws.Range("N135").Value = -50061

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 25466470
$ws.Range("K32").Value = 25466470
$ws.Range("M32").Value = -25466183
$ws.Range("H61").Value = 2720.0278
$ws.Range("I61").Value = 2219.182
$ws.Range("J61").Value = 3507.0715
$ws.Range("K61").Value = 2219.182
$ws.Range("L61").Value = 3507.0715
$ws.Range("M61").Value = -2007.182
$ws.Range("N61").Value = -3931.0715
$ws.Range("H63").Value = 4285.3335
$ws.Range("I63").Value = 2166.6667
$ws.Range("J63").Value = 4815
$ws.Range("K63").Value = 2166.6667
$ws.Range("L63").Value = 4815
$ws.Range("M63").Value = -1480.6667
$ws.Range("N63").Value = -6187
$ws.Range("H66").Value = 4285.3335
$ws.Range("I66").Value = 2166.6667
$ws.Range("J66").Value = 4815
$ws.Range("K66").Value = 10833.3335
$ws.Range("L66").Value = 24075
$ws.Range("M66").Value = -7401.333500000001
$ws.Range("N66").Value = -30939
$ws.Range("H122").Value = 3803.8
$ws.Range("I122").Value = 3378.182
$ws.Range("J122").Value = 4974.25
$ws.Range("K122").Value = 10134.546
$ws.Range("L122").Value = 14922.75
$ws.Range("M122").Value = -7684.545999999998
$ws.Range("N122").Value = -19822.75
$ws.Range("H132").Value = 316595.97
$ws.Range("I132").Value = 479145.84
$ws.Range("J132").Value = 6273.4546
$ws.Range("K132").Value = 1437437.52
$ws.Range("L132").Value = 18820.3638
$ws.Range("M132").Value = -1434907.52
$ws.Range("N132").Value = -23880.3638
$ws.Range("H136").Value = 2720.0278
$ws.Range("I136").Value = 2219.182
$ws.Range("J136").Value = 3507.0715
$ws.Range("K136").Value = 6657.545999999999
$ws.Range("L136").Value = 10521.2145
$ws.Range("M136").Value = -4107.545999999999
$ws.Range("N136").Value = -15621.2145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2466183.2
$ws.Range("I134").Value = 3107756.5
$ws.Range("J134").Value = 6818.8335
$ws.Range("K134").Value = 9323269.5
$ws.Range("L134").Value = 20456.5005
$ws.Range("M134").Value = -9320734.5
$ws.Range("N134").Value = -25526.5005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4350.5884
$ws.Range("I31").Value = 2062.0667
$ws.Range("J31").Value = 6157.316
$ws.Range("K31").Value = 2062.0667
$ws.Range("L31").Value = 6157.316
$ws.Range("M31").Value = -1767.0667
$ws.Range("N31").Value = -6747.316
$ws.Range("H34").Value = 4350.5884
$ws.Range("I34").Value = 2062.0667
$ws.Range("J34").Value = 6157.316
$ws.Range("K34").Value = 2062.0667
$ws.Range("L34").Value = 6157.316
$ws.Range("M34").Value = -1860.0667
$ws.Range("N34").Value = -6561.316
$ws.Range("H58").Value = 2362.276
$ws.Range("I58").Value = 2107.617
$ws.Range("K58").Value = 2107.617
$ws.Range("M58").Value = -1904.617
$ws.Range("H87").Value = 106999
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 106999
$ws.Range("K87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("M87").Value = 106999
$ws.Range("N87").Value = -109371
$ws.Range("H90").Value = 106999
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 106999
$ws.Range("K90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("M90").Value = 320997
$ws.Range("N90").Value = -332853
$ws.Range("H136").Value = 2362.276
$ws.Range("I136").Value = 2107.617
$ws.Range("K136").Value = 6322.851000000001
$ws.Range("M136").Value = -3772.851000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2002.2142
$ws.Range("I5").Value = 1226.4445
$ws.Range("J5").Value = 3398.6
$ws.Range("K5").Value = 3679.3335
$ws.Range("L5").Value = 10195.8
$ws.Range("M5").Value = -3567.3335
$ws.Range("N5").Value = -10419.8
$ws.Range("H135").Value = 2002.2142
$ws.Range("I135").Value = 1226.4445
$ws.Range("J135").Value = 3398.6
$ws.Range("K135").Value = 11038.0005
$ws.Range("L135").Value = 30587.4
$ws.Range("M135").Value = -8503.0005
$ws.Range("N135").Value = -35657.39999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 13402.066
$ws.Range("I2").Value = 79.30768999999999
$ws.Range("K2").Value = 79.30768999999999
$ws.Range("M2").Value = 33.69231000000001
$ws.Range("H97").Value = 897.1053000000001
$ws.Range("I97").Value = 724.7222
$ws.Range("K97").Value = 724.7222
$ws.Range("M97").Value = -228.7222
$ws.Range("H122").Value = 1323.9333
$ws.Range("I122").Value = 927.6667
$ws.Range("K122").Value = 2783.0001
$ws.Range("M122").Value = -333.0001000000002

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2264.4333
$ws.Range("I82").Value = 2406.45
$ws.Range("J82").Value = 1980.4
$ws.Range("K82").Value = 2406.45
$ws.Range("L82").Value = 1980.4
$ws.Range("M82").Value = -2045.45
$ws.Range("N82").Value = -2702.4
$ws.Range("H85").Value = 2264.4333
$ws.Range("I85").Value = 2406.45
$ws.Range("J85").Value = 1980.4
$ws.Range("K85").Value = 2406.45
$ws.Range("L85").Value = 1980.4
$ws.Range("M85").Value = -1158.45
$ws.Range("N85").Value = -4476.4
$ws.Range("H132").Value = 1004088.5
$ws.Range("I132").Value = 1253249
$ws.Range("K132").Value = 3759747
$ws.Range("M132").Value = -3757217
$ws.Range("H136").Value = 25000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5900
$ws.Range("I62").Value = 4625
$ws.Range("J62").Value = 6466.6665
$ws.Range("K62").Value = 4625
$ws.Range("L62").Value = 6466.6665
$ws.Range("M62").Value = -4001
$ws.Range("N62").Value = -7714.6665
$ws.Range("H64").Value = 92487
$ws.Range("H65").Value = 5900
$ws.Range("I65").Value = 4625
$ws.Range("J65").Value = 6466.6665
$ws.Range("K65").Value = 23125
$ws.Range("L65").Value = 32333.3325
$ws.Range("M65").Value = -20005
$ws.Range("N65").Value = -38573.3325
$ws.Range("H67").Value = 92487
$ws.Range("H126").Value = 8087.25
$ws.Range("I126").Value = 8739.799999999999
$ws.Range("K126").Value = 26219.4
$ws.Range("M126").Value = -23749.4
$ws.Range("H132").Value = 39287.926
$ws.Range("I132").Value = 40645.152
$ws.Range("K132").Value = 121935.456
$ws.Range("M132").Value = -119405.456
